# Daily attendance processing - 2026-01-08 15:09:07
#
# For every row in the "Recorded By" column (G), when the cell lists
# multiple recorder names/emails separated by ", ", rotate the list one
# position to the left (the first entry moves to the end). Cells with a
# single entry are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
            $cell.Value = $rotated
        }
    }
}
